$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.215.85"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.654.74"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'219.11"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'0.5243"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D8").Value = "'0.2663"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'0.06358"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'20.69"
$ws.Range("D11").Value = "'0.07723"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "'4.597"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "1.606.93"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "1.884.08"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.5624"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "0.0₅8249"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'65.38"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "26.220.66"
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "'4.695"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'10.40"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "'191.98"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "'6.003"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "'143.68"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "'0.1204"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "'7.270"
$ws.Range("D28").Value = "'15.95"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "'1.515"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'3.505"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'3.359"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'0.9539"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "'2.801"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "'2.414"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "'0.5758"
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "'6.007"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'0.8416"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "'101.95"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "1.008.22"
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("D45").Value = "1.794.99"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'58.39"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'0.05344"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.039"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.4349"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₈102"
$ws.Range("E51").Value = "  -4.19%  "
